$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("test_add_valid_employee") ---
# Header row unchanged except Full name/Email columns reference new shared strings
$ws1.Range("C1").Value = "Full name"
$ws1.Range("D1").Value = "Email"
$ws1.Range("E1").Value = "Jobtitle"

# Data row: replace old candidate (padmakshi / princy / Ken kevin / kevin@gmai.com) with new bot data
$ws1.Range("B2").Value = "bot123"
$ws1.Range("A2").Value = "rohit@beheraemail.com"
$ws1.Range("C2").Value = "Bot Rohit"
$ws1.Range("D2").Value = "bot@behera.com"
$ws1.Range("E2").Value = "Engineer"

# Remove hyperlinks that existed on A2 and D2 (keep formatting/style)
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Range("D2").Hyperlinks.Delete()

# --- Sheet2 ("test_invalid_profile_upload") ---
$ws2.Range("A1").Value = "Username"
$ws2.Range("C1").Value = "upload number"
$ws2.Range("D1").Value = "Expected Error"

$ws2.Range("A2").Value = "rohit@beheraemail.com"
$ws2.Range("B2").Value = "bot123"
$ws2.Range("C2").Value = 91767
$ws2.Range("D2").Value = "Please enter a valid phone number."

$ws2.Range("A2").Hyperlinks.Delete()

# --- Sheet view / selection changes ---
$ws1.Activate()
$ws1.Range("H12").Select()

# Make sheet2 the active sheet/tab and set its selection
$ws2.Activate()
$ws2.Range("D15").Select()
